$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new blank rows starting at row 12 (pushes the K:L lookup
#     tables that previously started at row 12 down to row 16, and the
#     "TtblTipoSeguro" block down to rows 21-26) ---
$ws.Rows("12:15").Insert()

# --- Two new policy records added to the tbl_Apolice table (A12:I13) ---
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4
$ws.Range("H12").Value = 125
$ws.Range("I12").Value = 30000

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 12
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4
$ws.Range("H13").Value = 650
$ws.Range("I13").Value = 45000

# --- Correct an existing data value: policy #5 (row 8) idSeguradora 4 -> 3 ---
$ws.Range("C8").Value = 3

# --- Move the "TBL Seguradora" mini-table header + column titles from
#     K10:L11 down into the blank rows opened up at K14:L15, then clear
#     the vacated source cells ---
$ws.Range("K14").Value = $ws.Range("K10").Text
$ws.Range("K15").Value = $ws.Range("K11").Text
$ws.Range("L15").Value = $ws.Range("L11").Text

$ws.Range("K10").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()

# --- Rename the last "TtblTipoSeguro" lookup entry from "Carro" to
#     "Automóvel" (now at L26 after the row insert) ---
$ws.Range("L26").Value = "Automóvel"

# --- Update the active selection to reflect the author's final cursor
#     position ---
$ws.Range("D5").Select()
